$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "RequestProcessingType" (col F) and "security" (col N) fields are no
# longer part of the standard field set, so drop both columns outright.
# Everything to their right shifts left to close the gap (G..P -> F..N).
$ws.Range("N1").EntireColumn.Delete()
$ws.Range("F1").EntireColumn.Delete()

# After the shift, rename a few headers to the new standard field names
# (HTTPAction -> Action, ExcludeField -> ExcludeFields, HttpStatusCode ->
# StatusCode) and capitalize the former "tags" header, which has landed in
# column M, to "Tags".
$ws.Range("H1").Value = "Action"
$ws.Range("I1").Value = "ExcludeFields"
$ws.Range("J1").Value = "StatusCode"
$ws.Range("M1").Value = "Tags"

# The "id;name;category.id" rich-text cell (now N3, was P3) previously had
# no explicit run formatting; give both runs an explicit Calibri 12 face,
# keeping the existing blue/underlined styling on "category.id".
$c = $ws.Range("N3")
$c.Characters(1, 8).Font.Name = "Calibri"
$c.Characters(1, 8).Font.Size = 12
$c.Characters(9, 12).Font.Name = "Calibri"
$c.Characters(9, 12).Font.Size = 12
$c.Characters(9, 12).Font.Underline = $true

# Match the saved selection/active cell.
$ws.Range("M1").Select()
